$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($Sheet, $Row, $Col, $Text) {
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# Row 2
Set-CellText $ws 2 4 '26.684.33'
Set-CellText $ws 2 5 '  +1.83%  '

# Row 3
Set-CellText $ws 3 4 '1.635.63'
Set-CellText $ws 3 5 '  +2.03%  '

# Row 4
Set-CellText $ws 4 5 '  -0.03%  '

# Row 5
Set-CellText $ws 5 4 '212.60'
Set-CellText $ws 5 5 '  +0.13%  '

# Row 6
Set-CellText $ws 6 4 '0.494'
Set-CellText $ws 6 5 '  +2.26%  '

# Row 7
Set-CellText $ws 7 5 '  -0.06%  '

# Row 8
Set-CellText $ws 8 4 '0.252'
Set-CellText $ws 8 5 '  +1.46%  '

# Row 9
Set-CellText $ws 9 5 '  +1.83%  '

# Row 10
Set-CellText $ws 10 4 '19.01'
Set-CellText $ws 10 5 '  +4.14%  '

# Row 11
Set-CellText $ws 11 5 '  +2.49%  '

# Row 12
Set-CellText $ws 12 4 '1.862.30'
Set-CellText $ws 12 5 '  +1.94%  '

# Row 13
Set-CellText $ws 13 4 '1.627.97'
Set-CellText $ws 13 5 '  +1.43%  '

# Row 14
Set-CellText $ws 14 4 '4.07'
Set-CellText $ws 14 5 '  +1.52%  '

# Row 15
Set-CellText $ws 15 4 '0.525'
Set-CellText $ws 15 5 '  +2.82%  '

# Row 16
Set-CellText $ws 16 4 '26.677.04'
Set-CellText $ws 16 5 '  +1.95%  '

# Row 17
Set-CellText $ws 17 4 '63.01'
Set-CellText $ws 17 5 '  +1.99%  '

# Row 18
Set-CellText $ws 18 4 '0.0₃0740'

# Row 20
Set-CellText $ws 20 4 '208.62'
Set-CellText $ws 20 5 '  +4.18%  '

# Row 21
Set-CellText $ws 21 4 '4.31'
Set-CellText $ws 21 5 '  +0.90%  '

# Row 22
Set-CellText $ws 22 2 'Chainlink'
Set-CellText $ws 22 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText $ws 22 4 '6.17'
Set-CellText $ws 22 5 '  +3.19%  '

# Row 23
Set-CellText $ws 23 2 'Avalanche'
Set-CellText $ws 23 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-CellText $ws 23 4 '9.38'
Set-CellText $ws 23 5 '  +1.32%  '

# Row 24
Set-CellText $ws 24 4 '1.92'
Set-CellText $ws 24 5 '  +3.56%  '

# Row 25
Set-CellText $ws 25 4 '146.50'
Set-CellText $ws 25 5 '  +1.63%  '

# Row 26
Set-CellText $ws 26 5 '  +0.03%  '

# Row 27
Set-CellText $ws 27 5 '  -0.58%  '

# Row 28
Set-CellText $ws 28 4 '6.75'
Set-CellText $ws 28 5 '  +3.05%  '

# Row 29
Set-CellText $ws 29 4 '15.40'
Set-CellText $ws 29 5 '  +1.57%  '

# Row 30
Set-CellText $ws 30 4 '0.0517'
Set-CellText $ws 30 5 '  +5.99%  '

# Row 31
Set-CellText $ws 31 5 '  -0.20%  '

# Row 32
Set-CellText $ws 32 5 '  +1.47%  '

# Row 33
Set-CellText $ws 33 5 '  +1.19%  '

# Row 34
Set-CellText $ws 34 5 '  +1.26%  '

# Row 35
Set-CellText $ws 35 4 '1.50'
Set-CellText $ws 35 5 '  +0.90%  '

# Row 36
Set-CellText $ws 36 4 '1.168.20'
Set-CellText $ws 36 5 '  +0.62%  '

# Row 37
Set-CellText $ws 37 5 '  -0.63%  '

# Row 38
Set-CellText $ws 38 4 '0.808'
Set-CellText $ws 38 5 '  +2.84%  '

# Row 39
Set-CellText $ws 39 5 '  +0.04%  '

# Row 40
Set-CellText $ws 40 2 'MXToken'
Set-CellText $ws 40 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws 40 4 '2.32'
Set-CellText $ws 40 5 '  +0.43%  '

# Row 41
Set-CellText $ws 41 2 'ImmutableX'
Set-CellText $ws 41 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 41 4 '0.503'
Set-CellText $ws 41 5 '  +1.85%  '

# Row 42
Set-CellText $ws 42 4 '0.794'
Set-CellText $ws 42 5 '  +1.82%  '

# Row 43
Set-CellText $ws 43 5 '  +1.39%  '

# Row 44
Set-CellText $ws 44 4 '1.772.80'
Set-CellText $ws 44 5 '  +1.95%  '

# Row 45
Set-CellText $ws 45 4 '92.43'
Set-CellText $ws 45 5 '  +1.10%  '

# Row 46
Set-CellText $ws 46 4 '1.55'
Set-CellText $ws 46 5 '  +2.24%  '

# Row 47
Set-CellText $ws 47 5 '  +4.32%  '

# Row 48
Set-CellText $ws 48 4 '54.77'
Set-CellText $ws 48 5 '  +1.39%  '

# Row 49
Set-CellText $ws 49 5 '  +1.47%  '

# Row 50
Set-CellText $ws 50 4 '0.410'
Set-CellText $ws 50 5 '  +0.69%  '

# Row 51
Set-CellText $ws 51 5 '  +4.30%  '
